# Add a new worksheet "2016-10-11" right after the existing "2016-10-07"
# sheet. The new sheet is produced as a copy of the existing sheet (so it
# inherits the same sheet/page setup), then trimmed down to just the
# header row, and finally renamed.

$wb = $excel.ActiveWorkbook

$existing = $wb.Worksheets.Item("2016-10-07")

# Copy the existing sheet, placing the new copy immediately after it.
$existing.Copy($null, $existing)
$new = $wb.Worksheets.Item(2)
$new.Name = "2016-10-11"

# Keep only the header row (row 1); drop the data rows that came along
# with the copy so the new sheet starts out with just the column headers.
$new.Rows("2:51").Delete()

# Restore the original sheet as the active/selected tab.
$existing.Activate()
